$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1338
$ws1.Range("F8").Value = 11583
$ws1.Range("F9").Value = 4366
$ws1.Range("F14").Value = 2541
$ws1.Range("F16").Value = 138
$ws1.Range("F18").Value = 4496
$ws1.Range("F21").Value = 11323
$ws1.Range("F22").Value = 11243
$ws1.Range("F27").Value = 43

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 1338
$ws4.Range("F8").Value = 11583
$ws4.Range("F9").Value = 4366
$ws4.Range("F14").Value = 2541
$ws4.Range("F17").Value = 138
$ws4.Range("F19").Value = 4496
$ws4.Range("F22").Value = 11323
$ws4.Range("F23").Value = 11243
$ws4.Range("F28").Value = 43
